$d = $word.ActiveDocument

$pairs = @(
    @("2023-12-18 Monday", "2023-12-19 Tuesday"),
    @("39-35=4", "55-4=51"),
    @("47-13=34", "98-35=63"),
    @("75-67=8", "48+9=57"),
    @("59+26=85", "98-67=31"),
    @("93-65=28", "90-33=57"),
    @("19+43=62", "15+28=43"),
    @("87+5=92", "87-34=53"),
    @("42-30=12", "91-76=15"),
    @("24+18=42", "64+1=65"),
    @("88-64=24", "4+2=6"),
    @("82+5=87", "68-37=31"),
    @("4+50=54", "86-78=8"),
    @("74+21=95", "33+28=61"),
    @("58-11=47", "50+1=51"),
    @("56-14=42", "98-27=71"),
    @("54-52=2", "12+38=50"),
    @("27+6=33", "3+32=35"),
    @("60-14=46", "14+8=22"),
    @("14+55=69", "60+8=68"),
    @("8-8=0", "34+56=90"),
    @("76-54=22", "2+86=88"),
    @("91-86=5", "8+77=85"),
    @("31+15=46", "76-9=67"),
    @("9+43=52", "1+1=2"),
    @("57-34=23", "98-46=52"),
    @("80-6=74", "5+3=8"),
    @("47-39=8", "76+7=83"),
    @("21+43=64", "93-18=75"),
    @("86-46=40", "83-35=48"),
    @("11+88=99", "54-10=44"),
    @("11+16=27", "60-34=26"),
    @("90-37=53", "39-19=20"),
    @("84-41=43", "8+65=73"),
    @("92-47=45", "17+2=19"),
    @("66-39=27", "64-45=19"),
    @("96-44=52", "0+31=31"),
    @("16+27=43", "76-75=1"),
    @("71-27=44", "22+43=65"),
    @("29+27=56", "72-9=63"),
    @("1+66=67", "24+9=33"),
    @("24+25=49", "89-34=55"),
    @("74-59=15", "58-45=13"),
    @("44-1=43", "68+25=93"),
    @("28+24=52", "92-3=89"),
    @("82-7=75", "36-18=18"),
    @("76-14=62", "32-19=13"),
    @("57+8=65", "93-80=13"),
    @("70-66=4", "5+68=73"),
    @("10+8=18", "95-78=17"),
    @("60-2=58", "31+46=77"),
    @("56-11=45", "22+67=89"),
    @("23+24=47", "46+34=80"),
    @("80-48=32", "2+22=24"),
    @("92-41=51", "84-9=75"),
    @("81+8=89", "12+48=60"),
    @("24-20=4", "7+36=43"),
    @("68-54=14", "51-15=36"),
    @("50-4=46", "54-48=6"),
    @("37+62=99", "50-42=8"),
    @("46-36=10", "42+6=48"),
    @("64-30=34", "82-34=48"),
    @("8+52=60", "5+33=38"),
    @("11+69=80", "9+30=39"),
    @("6+3=9", "34-4=30"),
    @("15+70=85", "24+69=93"),
    @("65+14=79", "68-15=53"),
    @("74-66=8", "25+11=36"),
    @("94-59=35", "24+33=57"),
    @("35-16=19", "46-24=22"),
    @("59-6=53", "47+45=92"),
    @("14+10=24", "52+23=75"),
    @("92-48=44", "28-21=7"),
    @("70-42=28", "90-29=61"),
    @("37-14=23", "23+66=89"),
    @("60-45=15", "60+34=94"),
    @("51+27=78", "17+1=18"),
    @("30+37=67", "20+79=99"),
    @("9+12=21", "88-32=56"),
    @("19+4=23", "31+19=50"),
    @("12+10=22", "46+17=63"),
    @("58-51=7", "29+40=69"),
    @("75-20=55", "90-44=46"),
    @("49+46=95", "98-9=89"),
    @("57-48=9", "83-37=46"),
    @("22+23=45", "58-45=13"),
    @("14+35=49", "77-24=53"),
    @("43+12=55", "85-54=31"),
    @("98-1=97", "91-70=21"),
    @("47+9=56", "57+0=57"),
    @("58+14=72", "38-37=1"),
    @("89-13=76", "81-9=72"),
    @("91-5=86", "47+46=93"),
    @("5+4=9", "17+38=55"),
    @("13+44=57", "13+28=41"),
    @("64-16=48", "95-37=58"),
    @("48-46=2", "92-71=21"),
    @("84+8=92", "28+36=64"),
    @("69-49=20", "55-26=29"),
    @("91-9=82", "75-17=58"),
    @("97-82=15", "27-22=5"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
